# Atualizada documentação de projeto, validadores
#
# 1. Bump the version line from "Versão: 2.0" to "Versão: 3.0", with the
#    number split across its own runs (as Word does when you edit just the
#    "2" character in place).
# 2. Add a new "Revisado: 12/04/2025" paragraph right after the version line.

$d = $word.ActiveDocument

# --- Locate the "Versão: 2.0" paragraph -------------------------------------------------
$versionPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Versão: 2.0*") {
        $versionPara = $para
        break
    }
}

$r = $versionPara.Range

# The paragraph text is "Versão: 2.0" followed by the paragraph mark, so the
# last 4 characters before the mark are "2.0".
$numRange = $d.Range($r.End - 4, $r.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
    '<w:r>' + `
      '<w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:lang w:val="pt-BR"/></w:rPr>' + `
      '<w:t>3</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
      '<w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:lang w:val="pt-BR"/></w:rPr>' + `
      '<w:t>.0</w:t>' + `
    '</w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>'
$numRange.InsertXML($xml)

# --- Insert the new "Revisado: 12/04/2025" paragraph right after it ---------------------
$versionPara.Range.InsertParagraphAfter()
$revisadoPara = $versionPara.Next(1, 1)
$revisadoPara.Range.Text = "Revisado: 12/04/2025"
